# Pacientes.xlsx fix: add "rehabilitado" column, fill in patient data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$table = $ws.ListObjects.Item("Tabla1")

# Insert a new "rehabilitado" table column right before "tecnico"
# (tecnico is currently the 10th table column).
$tecnicoCol = $table.ListColumns.Item("tecnico")
$newCol = $table.ListColumns.Add($tecnicoCol.Index)
$newCol.Name = "rehabilitado"

# Comment explaining the new column's accepted values.
$ws.Range("J1").AddComment("Puede tener el valor ""X"" para marcar que esta rehabilitado, si se deja vacio se guardara como no rehabilitado")

# --- Update existing data row (row 2) ---
$ws.Range("A2").Value = "joaquin"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = "99186787T"
$ws.Range("G2").Value = "san marques"
$ws.Range("J2").Value = "X"
$ws.Range("L2").Value = "asdf"

# --- Add new data row (row 3) ---
$ws.Range("A3").Value = "maria"
$ws.Range("B3").Value = "cabrera"
$ws.Range("D3").Value = "65963475G"
$ws.Range("E3").Value = Get-Date -Year 1999 -Month 1 -Day 10
$ws.Range("F3").Value = "Mujer"
$ws.Range("I3").Value = "asdf"
$ws.Range("K3").Value = "jose maria"
$ws.Range("L3").Value = "nada"

$ws.Range("L3").Select()
